$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.216.24'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '1.859.94'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'236.01"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('D6').Value = "'1.001"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').Value = "'0.4717"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.09%  '
$ws.Range('D8').Value = "'0.2896"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.40%  '
$ws.Range('D9').Value = "'0.06561"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.33%  '
$ws.Range('D10').Value = "'21.81"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.99%  '
$ws.Range('E11').Value = '  +1.12%  '
$ws.Range('D12').Value = "'97.75"
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Value = '1.856.98'
$ws.Range('E13').Value = '  -0.46%  '
$ws.Range('D14').Value = "'5.130"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.58%  '
$ws.Range('D15').Value = "'0.6803"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.08%  '
$ws.Range('D16').Value = "'266.64"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.62%  '
$ws.Range('D17').Value = '30.216.68'
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('D18').Value = "'13.67"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +8.11%  '
$ws.Range('D19').Value = "'0.000007545"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.54%  '
$ws.Range('E20').Value = '  -0.08%  '
$ws.Range('D21').Value = '2.097.94'
$ws.Range('E21').Value = '  -0.77%  '
$ws.Range('D22').Value = "'1.001"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').Value = "'5.262"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.79%  '
$ws.Range('D24').Value = "'6.167"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').Value = "'167.50"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.43%  '
$ws.Range('D26').Value = "'9.182"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').Value = "'18.91"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.07%  '
$ws.Range('D28').Value = "'1.948"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.46%  '
$ws.Range('D29').Value = "'1.394"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.27%  '
$ws.Range('D30').Value = "'0.09924"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.26%  '
$ws.Range('D31').Value = "'4.331"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.95%  '
$ws.Range('D32').Value = "'1.469"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('D33').Value = "'4.008"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.05%  '
$ws.Range('D34').Value = "'0.04706"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.39%  '
$ws.Range('D35').Value = "'1.129"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.06%  '
$ws.Range('D36').Value = "'0.7009"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.52%  '
$ws.Range('D37').Value = "'2.707"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.70%  '
$ws.Range('D38').Value = "'0.01876"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.53%  '
$ws.Range('D39').Value = "'2.620"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.34%  '
$ws.Range('D40').Value = "'6.324"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.04%  '
$ws.Range('D41').Value = "'73.85"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.40%  '
$ws.Range('D42').Value = "'1.940"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('D43').Value = "'0.8414"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.52%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = "'0.9999"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value = "'0.4156"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.22%  '
$ws.Range('D46').Value = "'103.26"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('D47').Value = "'7.139"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.48%  '
$ws.Range('D48').Value = "'942.54"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.78%  '
$ws.Range('D49').Value = "'9.212"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('D50').Value = "'34.12"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.46%  '
$ws.Range('E51').Value = '  +0.52%  '
